$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 107.62
$ws.Range("I15").Value = 107.62
$ws.Range("K15").Value = 322.86
$ws.Range("M15").Value = -153.86
$ws.Range("H40").Value = 1408
$ws.Range("I40").Value = 1380
$ws.Range("J40").Value = 1450
$ws.Range("K40").Value = 1380
$ws.Range("L40").Value = 1450
$ws.Range("M40").Value = -1205
$ws.Range("N40").Value = -1800
$ws.Range("H112").Value = 1488.75
$ws.Range("J112").Value = 1494.697
$ws.Range("L112").Value = 4484.090999999999
$ws.Range("N112").Value = -6700.090999999999
$ws.Range("H132").Value = 19309584
$ws.Range("I132").Value = 20491260
$ws.Range("K132").Value = 61473780
$ws.Range("M132").Value = -61471250

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8416.618
$ws.Range("I32").Value = 5272.9375
$ws.Range("K32").Value = 5272.9375
$ws.Range("M32").Value = -4985.9375
$ws.Range("H74").Value = 1098.2808
$ws.Range("I74").Value = 807.5306399999999
$ws.Range("J74").Value = 2879.125
$ws.Range("K74").Value = 807.5306399999999
$ws.Range("L74").Value = 2879.125
$ws.Range("M74").Value = 66.46936000000005
$ws.Range("N74").Value = -4627.125
$ws.Range("H77").Value = 1098.2808
$ws.Range("I77").Value = 807.5306399999999
$ws.Range("J77").Value = 2879.125
$ws.Range("K77").Value = 4037.6532
$ws.Range("L77").Value = 14395.625
$ws.Range("M77").Value = 330.3468000000003
$ws.Range("N77").Value = -23131.625
$ws.Range("H122").Value = 2311.923
$ws.Range("I122").Value = 1225.5
$ws.Range("K122").Value = 3676.5
$ws.Range("M122").Value = -1226.5
$ws.Range("H132").Value = 2445.9692
$ws.Range("I132").Value = 1864.6666
$ws.Range("J132").Value = 5299.636
$ws.Range("K132").Value = 5593.9998
$ws.Range("L132").Value = 15898.908
$ws.Range("M132").Value = -3063.9998
$ws.Range("N132").Value = -20958.908
$ws.Range("H137").Value = 50566.668
$ws.Range("J137").Value = 50566.668
$ws.Range("L137").Value = 50566.668
$ws.Range("N137").Value = -60766.668

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 894.44446
$ws.Range("I22").Value = 464.2857
$ws.Range("K22").Value = 464.2857
$ws.Range("M22").Value = -291.2857
$ws.Range("H99").Value = 3512.5652
$ws.Range("I99").Value = 1962.7273
$ws.Range("J99").Value = 4933.25
$ws.Range("K99").Value = 1962.7273
$ws.Range("L99").Value = 4933.25
$ws.Range("M99").Value = -464.7273
$ws.Range("N99").Value = -7929.25
$ws.Range("H134").Value = 2003.8372
$ws.Range("I134").Value = 1226.7894
$ws.Range("J134").Value = 7909.4
$ws.Range("K134").Value = 3680.3682
$ws.Range("L134").Value = 23728.2
$ws.Range("M134").Value = -1145.3682
$ws.Range("N134").Value = -28798.2
$ws.Range("H137").Value = 32962.5
$ws.Range("J137").Value = 32962.5
$ws.Range("L137").Value = 32962.5
$ws.Range("N137").Value = -43162.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2781.093
$ws.Range("I31").Value = 1004.7778
$ws.Range("J31").Value = 5778.625
$ws.Range("K31").Value = 1004.7778
$ws.Range("L31").Value = 5778.625
$ws.Range("M31").Value = -709.7778
$ws.Range("N31").Value = -6368.625
$ws.Range("H34").Value = 2781.093
$ws.Range("I34").Value = 1004.7778
$ws.Range("J34").Value = 5778.625
$ws.Range("K34").Value = 1004.7778
$ws.Range("L34").Value = 5778.625
$ws.Range("M34").Value = -802.7778
$ws.Range("N34").Value = -6182.625

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H40").Value = 459
$ws.Range("I40").Value = 150.3
$ws.Range("J40").Value = 900
$ws.Range("K40").Value = 601.2
$ws.Range("L40").Value = 3600
$ws.Range("M40").Value = -532.2
$ws.Range("N40").Value = -3738
$ws.Range("H62").Value = 7735
$ws.Range("J62").Value = 9980
$ws.Range("L62").Value = 29940
$ws.Range("N62").Value = -31312
$ws.Range("H65").Value = 7735
$ws.Range("J65").Value = 9980
$ws.Range("L65").Value = 89820
$ws.Range("N65").Value = -96684
$ws.Range("H81").Value = 1004.3
$ws.Range("I81").Value = 363.2857
$ws.Range("K81").Value = 1089.8571
$ws.Range("M81").Value = 33.14289999999983
$ws.Range("H84").Value = 1004.3
$ws.Range("I84").Value = 363.2857
$ws.Range("K84").Value = 3269.5713
$ws.Range("M84").Value = 2346.4287
$ws.Range("H122").Value = 2613.2666
$ws.Range("J122").Value = 3366.509
$ws.Range("L122").Value = 30298.581
$ws.Range("N122").Value = -35198.581
$ws.Range("H125").Value = 5800
$ws.Range("J125").Value = 7857.143
$ws.Range("L125").Value = 23571.429
$ws.Range("N125").Value = -33411.429
$ws.Range("H131").Value = 7813507
$ws.Range("J131").Value = 884.3771
$ws.Range("L131").Value = 2653.1313
$ws.Range("N131").Value = -12733.1313

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2718.0312
$ws.Range("I102").Value = 1981.6957
$ws.Range("J102").Value = 4599.778
$ws.Range("K102").Value = 1981.6957
$ws.Range("L102").Value = 4599.778
$ws.Range("M102").Value = -359.6957
$ws.Range("N102").Value = -7843.778
$ws.Range("H132").Value = 3063.4482
$ws.Range("I132").Value = 1960.8
$ws.Range("J132").Value = 5513.778
$ws.Range("K132").Value = 5882.4
$ws.Range("L132").Value = 16541.334
$ws.Range("M132").Value = -3352.4
$ws.Range("N132").Value = -21601.334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2007.2142
$ws.Range("I16").Value = 1830.1
$ws.Range("J16").Value = 2450
$ws.Range("K16").Value = 1830.1
$ws.Range("L16").Value = 2450
$ws.Range("M16").Value = -1660.1
$ws.Range("N16").Value = -2790
$ws.Range("H68").Value = 950.4761999999999
$ws.Range("I68").Value = 848
$ws.Range("J68").Value = 3000
$ws.Range("K68").Value = 848
$ws.Range("L68").Value = 3000
$ws.Range("M68").Value = -99
$ws.Range("N68").Value = -4498
$ws.Range("H71").Value = 950.4761999999999
$ws.Range("I71").Value = 848
$ws.Range("J71").Value = 3000
$ws.Range("K71").Value = 4240
$ws.Range("L71").Value = 15000
$ws.Range("M71").Value = -496
$ws.Range("N71").Value = -22488
$ws.Range("H122").Value = 4066.5715
$ws.Range("I122").Value = 3147.7
$ws.Range("J122").Value = 9579.799999999999
$ws.Range("K122").Value = 9443.099999999999
$ws.Range("L122").Value = 28739.4
$ws.Range("M122").Value = -6993.099999999999
$ws.Range("N122").Value = -33639.39999999999
$ws.Range("H132").Value = 5228.021
$ws.Range("I132").Value = 1649.8334
$ws.Range("K132").Value = 4949.5002
$ws.Range("M132").Value = -2419.5002

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 301214.28
$ws.Range("I126").Value = 2440.5908
$ws.Range("J126").Value = 770715.8
$ws.Range("K126").Value = 7321.7724
$ws.Range("L126").Value = 2312147.4
$ws.Range("M126").Value = -4851.7724
$ws.Range("N126").Value = -2317087.4
$ws.Range("H132").Value = 7938461.5
$ws.Range("I132").Value = 715.0625
$ws.Range("J132").Value = 12823228
$ws.Range("K132").Value = 2145.1875
$ws.Range("L132").Value = 38469684
$ws.Range("M132").Value = 384.8125
$ws.Range("N132").Value = -38474744
$ws.Range("H136").Value = 3138.2974
$ws.Range("I136").Value = 949.2381
$ws.Range("J136").Value = 6011.4375
$ws.Range("K136").Value = 2847.7143
$ws.Range("L136").Value = 18034.3125
$ws.Range("M136").Value = -297.7143000000001
$ws.Range("N136").Value = -23134.3125
